# Atualizado por script em 12-11-2023 20:45
# Adds two new match rows (90 and 91) to the end of the betting-odds sheet,
# mirroring the existing layout/formatting of the last data row (89).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting (styles, number formats, borders) of the last existing
# data row onto the two new rows so the new cells match the sheet's look.
$ws.Range("A89:V89").Copy($ws.Range("A90:V91"))

# ---- Row 90: Leixoes 0 - 1 Maritimo ----
$ws.Range("A90").Value = 89
$ws.Range("B90").Value = "portugal"
$ws.Range("C90").Value = "liga-portugal-2"
$ws.Range("D90").Value = "2023-2024"
$ws.Range("E90").Value = 45242.6875
$ws.Range("F90").Value = "Leixoes"
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = "Maritimo"
$ws.Range("I90").Value = 1
$ws.Range("J90").Value = 3.2
$ws.Range("K90").Value = "08/11/2023 06:12"
$ws.Range("L90").Value = 3.79
$ws.Range("M90").Value = "12/11/2023 16:20"
$ws.Range("N90").Value = 3.39
$ws.Range("O90").Value = "08/11/2023 06:12"
$ws.Range("P90").Value = 3.4
$ws.Range("Q90").Value = "12/11/2023 16:20"
$ws.Range("R90").Value = 2.29
$ws.Range("S90").Value = "08/11/2023 06:12"
$ws.Range("T90").Value = 2.1
$ws.Range("U90").Value = "12/11/2023 16:20"
$ws.Range("V90").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/leixoes-maritimo/pWvUsUrb/"

# ---- Row 91: Academico Viseu 1 - 0 Benfica B ----
$ws.Range("A91").Value = 90
$ws.Range("B91").Value = "portugal"
$ws.Range("C91").Value = "liga-portugal-2"
$ws.Range("D91").Value = "2023-2024"
$ws.Range("E91").Value = 45242.79166666666
$ws.Range("F91").Value = "Academico Viseu"
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = "Benfica B"
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2.07
$ws.Range("K91").Value = "06/11/2023 19:12"
$ws.Range("L91").Value = 1.94
$ws.Range("M91").Value = "12/11/2023 18:35"
$ws.Range("N91").Value = 3.58
$ws.Range("O91").Value = "06/11/2023 19:12"
$ws.Range("P91").Value = 3.68
$ws.Range("Q91").Value = "12/11/2023 18:35"
$ws.Range("R91").Value = 3.34
$ws.Range("S91").Value = "06/11/2023 19:12"
$ws.Range("T91").Value = 4.01
$ws.Range("U91").Value = "12/11/2023 18:35"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/academico-viseu-benfica/neyovnrN/"

$wb.Save()
